# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Wed Nov 22 22:46:49 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E hold prices/percentages as plain text (e.g. "37.503.92", "1.00",
# "0.0₃0812"). Whenever the new text still looks like a number, force the cell
# to Text first so COM does not silently coerce it (e.g. "1.00" -> 1, dropping the
# trailing zero), then restore the default "Normal" style so no formatting leaks in.
function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $looksNumeric = $value -match "^-?[0-9]+(\.[0-9]+)?$"
    if ($looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}

# Row 2
Set-TextValue "D2" "37.503.92"
Set-TextValue "E2" "  +2.76%  "

# Row 3
Set-TextValue "D3" "2.070.07"
Set-TextValue "E3" "  +5.74%  "

# Row 4
Set-TextValue "E4" "  -0.28%  "

# Row 5
Set-TextValue "D5" "236.61"
Set-TextValue "E5" "  +3.32%  "

# Row 6
Set-TextValue "D6" "0.618"
Set-TextValue "E6" "  +3.68%  "

# Row 7
Set-TextValue "D7" "58.06"
Set-TextValue "E7" "  +9.85%  "

# Row 8
Set-TextValue "E8" "  -0.05%  "

# Row 9
Set-TextValue "E9" "  +4.51%  "

# Row 10
Set-TextValue "D10" "57.70"
Set-TextValue "E10" "  +0.91%  "

# Row 11
Set-TextValue "D11" "0.0762"
Set-TextValue "E11" "  +3.88%  "

# Row 12
Set-TextValue "E12" "  +4.40%  "

# Row 13
Set-TextValue "D13" "2.373.59"
Set-TextValue "E13" "  +5.63%  "

# Row 14
Set-TextValue "D14" "14.31"
Set-TextValue "E14" "  +4.16%  "

# Row 15
Set-TextValue "D15" "20.97"
Set-TextValue "E15" "  +7.00%  "

# Row 16
Set-TextValue "D16" "0.778"
Set-TextValue "E16" "  +5.20%  "

# Row 17
Set-TextValue "D17" "5.19"
Set-TextValue "E17" "  +4.61%  "

# Row 18
Set-TextValue "D18" "2.069.14"
Set-TextValue "E18" "  +5.45%  "

# Row 19
Set-TextValue "D19" "37.653.75"
Set-TextValue "E19" "  +3.17%  "

# Row 20
Set-TextValue "D20" "6.11"
Set-TextValue "E20" "  +23.25%  "

# Row 21
Set-TextValue "D21" "68.59"
Set-TextValue "E21" "  +2.15%  "

# Row 22
Set-TextValue "D22" "0.0₃0812"
Set-TextValue "E22" "  +2.74%  "

# Row 23
Set-TextValue "D23" "225.01"
Set-TextValue "E23" "  +2.38%  "

# Row 24
Set-TextValue "D24" "1.00"

# Row 25
Set-TextValue "B25" "PancakeSwap"
Set-TextValue "C25" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D25" "2.44"
Set-TextValue "E25" "  +5.74%  "

# Row 26
Set-TextValue "B26" "Toncoin"
Set-TextValue "C26" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D26" "2.40"
Set-TextValue "E26" "  +3.24%  "

# Row 27
Set-TextValue "D27" "162.91"
Set-TextValue "E27" "  +1.60%  "

# Row 28
Set-TextValue "E28" "  +4.46%  "

# Row 29
Set-TextValue "D29" "0.131"
Set-TextValue "E29" "  +8.02%  "

# Row 30
Set-TextValue "B30" "EthereumClassic"
Set-TextValue "C30" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D30" "19.35"
Set-TextValue "E30" "  +2.84%  "

# Row 31
Set-TextValue "B31" "ImmutableX"
Set-TextValue "C31" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D31" "1.39"
Set-TextValue "E31" "  +7.60%  "

# Row 32
Set-TextValue "D32" "0.119"
Set-TextValue "E32" "  +2.59%  "

# Row 33
Set-TextValue "B33" "LidoDAOToken"
Set-TextValue "C33" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D33" "2.63"
Set-TextValue "E33" "  +16.84%  "

# Row 34
Set-TextValue "B34" "Hedera"
Set-TextValue "C34" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D34" "0.0630"
Set-TextValue "E34" "  +5.15%  "

# Row 35
Set-TextValue "D35" "4.48"
Set-TextValue "E35" "  +3.98%  "

# Row 36
Set-TextValue "D36" "4.47"
Set-TextValue "E36" "  +7.49%  "

# Row 37
Set-TextValue "E37" "  -0.19%  "

# Row 38
Set-TextValue "D38" "1.79"
Set-TextValue "E38" "  +0.61%  "

# Row 39
Set-TextValue "D39" "3.36"
Set-TextValue "E39" "  +5.91%  "

# Row 40
Set-TextValue "D40" "5.86"
Set-TextValue "E40" "  +14.58%  "

# Row 41
Set-TextValue "E41" "  -1.08%  "

# Row 42
Set-TextValue "B42" "FTXToken"
Set-TextValue "C42" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D42" "4.44"
Set-TextValue "E42" "  +30.61%  "

# Row 43
Set-TextValue "B43" "Cronos"
Set-TextValue "C43" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D43" "0.0954"
Set-TextValue "E43" "  +10.10%  "

# Row 44
Set-TextValue "B44" "Maker"
Set-TextValue "C44" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D44" "1.475.54"
Set-TextValue "E44" "  +5.03%  "

# Row 45
Set-TextValue "B45" "Aave"
Set-TextValue "C45" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D45" "95.76"
Set-TextValue "E45" "  +11.06%  "

# Row 46
Set-TextValue "E46" "  +5.77%  "

# Row 47
Set-TextValue "D47" "16.15"
Set-TextValue "E47" "  +9.82%  "

# Row 48
Set-TextValue "E48" "  +4.41%  "

# Row 49
Set-TextValue "D49" "7.30"
Set-TextValue "E49" "  +9.65%  "

# Row 50
Set-TextValue "E50" "  +4.38%  "

# Row 51
Set-TextValue "D51" "2.93"
